$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking values in the Price (D) and Hora (G)
# columns so Excel does not auto-convert these text cells into numbers.
# Only touch the specific cells that receive new numeric-looking text so we do not
# disturb the number format of any untouched cell.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range('D2').Value = '247.84'
$ws.Range('G2').Value = '12'
$ws.Range('D3').Value = '21.69'
$ws.Range('G3').Value = '12'
$ws.Range('D4').Value = '5.426'
$ws.Range('G4').Value = '12'
$ws.Range('G5').Value = '12'
$ws.Range('D6').Value = '3.384'
$ws.Range('G6').Value = '12'
$ws.Range('D7').Value = '0.8076'
$ws.Range('G7').Value = '12'
$ws.Range('G8').Value = '12'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').Value = '0.1490'
$ws.Range('E9').Value = '8WazirXWRX'
$ws.Range('G9').Value = '12'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').Value = '0.08133'
$ws.Range('E10').Value = '9MandalaExchangeTokenMDX'
$ws.Range('G10').Value = '12'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '0.03159'
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('G11').Value = '12'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.03027'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('G12').Value = '12'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.09284'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('G13').Value = '12'
$ws.Range('B14').Value = 'MCDex'
$ws.Range('C14').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D14').Value = '3.527'
$ws.Range('E14').Value = '13MCDexMCB'
$ws.Range('G14').Value = '12'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001647'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('G15').Value = '12'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').Value = '0.04704'
$ws.Range('E16').Value = '15CoinExTokenCET'
$ws.Range('G16').Value = '12'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').Value = '0.0005867'
$ws.Range('E17').Value = '16OneONE'
$ws.Range('G17').Value = '12'
$ws.Range('D18').Value = '0.006364'
$ws.Range('G18').Value = '12'
$ws.Range('D19').Value = '0.005033'
$ws.Range('E19').Value = '18HotbitTokenHTBBestin24h'
$ws.Range('G19').Value = '12'
$ws.Range('G20').Value = '12'
$ws.Range('G21').Value = '12'
$ws.Range('G22').Value = '12'
$ws.Range('D23').Value = '3.773'
$ws.Range('G23').Value = '12'
$ws.Range('D24').Value = '6.433'
$ws.Range('G24').Value = '12'
$ws.Range('D25').Value = '2.146'
$ws.Range('G25').Value = '12'
$ws.Range('G26').Value = '12'
$ws.Range('G27').Value = '12'
$ws.Range('G28').Value = '12'
$ws.Range('G29').Value = '12'
$ws.Range('G30').Value = '12'
$ws.Range('G31').Value = '12'
$ws.Range('G32').Value = '12'
$ws.Range('G33').Value = '12'
$ws.Range('G34').Value = '12'
$ws.Range('G35').Value = '12'
$ws.Range('G36').Value = '12'
$ws.Range('G37').Value = '12'
$ws.Range('G38').Value = '12'
$ws.Range('G39').Value = '12'
$ws.Range('D40').Value = '0.04107'
$ws.Range('G40').Value = '12'
$ws.Range('D41').Value = '0.006960'
$ws.Range('G41').Value = '12'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Value = '0.1045'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('G42').Value = '12'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').Value = '0.002974'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('G43').Value = '12'
$ws.Range('D44').Value = '0.008645'
$ws.Range('G44').Value = '12'
$ws.Range('D45').Value = '0.00005898'
$ws.Range('G45').Value = '12'
$ws.Range('G46').Value = '12'
$ws.Range('D47').Value = '0.0005506'
$ws.Range('G47').Value = '12'
$ws.Range('D48').Value = '0.6832'
$ws.Range('G48').Value = '12'
$ws.Range('D49').Value = '0.008519'
$ws.Range('G49').Value = '12'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('G50').Value = '12'
$ws.Range('G51').Value = '12'
